# Refresh the crypto price/volume snapshot on Sheet1 (GitHub Actions bot run).
# Columns: D = Price (text, may look numeric), E = Volume(1h) (padded percent text).
# A leading '' (single quote) forces Excel to keep single-dot values like
# "1.002" as literal text instead of auto-converting them to a number,
# matching the original inlineStr/text storage of these cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.883.97'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.633.82'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''215.86'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = '''0.5085'
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.2579'
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").Value = '''0.06338'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = '''19.47'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").Value = '''0.07771'
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("D12").Value = '''4.265'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.634.94'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '1.859.06'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").Value = '''0.5508'
$ws.Range("D16").Value = '''63.90'
$ws.Range("D17").Value = '0.0₅7662'
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = '25.911.52'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '''195.14'
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").Value = '''4.416'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '''9.905'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '''1.003'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").Value = '''1.916'
$ws.Range("E25").Value = '  +3.39%  '
$ws.Range("D26").Value = '''142.58'
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").Value = '''0.1252'
$ws.Range("E27").Value = '  +5.75%  '
$ws.Range("D28").Value = '''6.771'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '''15.56'
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = '''1.239'
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").Value = '''0.04904'
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '''3.250'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = '''3.196'
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").Value = '''2.369'
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").Value = '''0.8990'
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").Value = '''0.5534'
$ws.Range("E37").Value = '  +2.93%  '
$ws.Range("D38").Value = '''2.543'
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = '1.116.99'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").Value = '''0.01557'
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("D41").Value = '''1.001'
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '''5.622'
$ws.Range("E42").Value = '  +3.53%  '
$ws.Range("D43").Value = '''0.7973'
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").Value = '''97.44'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("E45").Value = '  -5.29%  '
$ws.Range("D46").Value = '1.770.80'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").Value = '''1.004'
$ws.Range("D49").Value = '''54.82'
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").Value = '''0.05134'
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("D51").Value = '''7.616'
$ws.Range("E51").Value = '  +4.56%  '
